# Update the timestamp portion of the test e-mail addresses that were
# generated for this test run (e.g. juan.perez+20251109_022039@test.com
# -> juan.perez+20251109_024842@test.com). The same e-mail text is
# reused (as the very same shared string) on both the "UsuariosRegistro"
# sheet (column C, rows 2-6) and the "LoginData" sheet (column A, the
# rows that mirror the two "Success" logins), so every cell containing
# one of the old addresses must be refreshed to keep them in sync.

$wb = $excel.ActiveWorkbook

$oldStamp = "20251109_022039"
$newStamp = "20251109_024842"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($row in 1..$used.Rows.Count) {
        foreach ($col in 1..$used.Columns.Count) {
            $cell = $used.Cells.Item($row, $col)
            $text = $cell.Value2
            if ($text -and $text.ToString().Contains($oldStamp)) {
                $cell.Value2 = $text -replace $oldStamp, $newStamp
            }
        }
    }
}
